# Minor fixes here and there.
$wb = $excel.ActiveWorkbook

$npcs = $wb.Worksheets.Item("NPCs")
$cmds = $wb.Worksheets.Item("Npcs Commands")

# Fix the "Dungeon Maiden" row on the NPCs sheet so the display name matches
# the internal/real name used elsewhere ("DungeonMaiden"), and correct the
# x/y position values for that NPC.
$npcs.Range("B8").Value = "DungeonMaiden"
$npcs.Range("H8").Value = 448
$npcs.Range("I8").Value = 96

# Keep the Npcs Commands sheet's npc_id reference consistent with the rename.
$cmds.Range("A8").Value = "DungeonMaiden"
